$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.169.30"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.869.03"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "696.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "3.865.99"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.06%  "
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "4.521.67"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "3.868.14"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "71.215.97"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "498.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.02"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.72"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "3.822.52"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("E39").Value = "  +8.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.44"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  +8.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.02"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.83"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000311"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.45"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "417.05"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.74"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("E51").Value = "  -2.12%  "
